$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- 1. Move the block of rows 912:1010 (99 rows) so it lands right before
#        row 626. This shifts the former rows 626:911 down by 99 rows.
#        (Equivalent to select rows 912:1010, cut, insert before row 626.)
$ws.Rows("626:724").Insert()
$src = $ws.Range("A1011:B1109")
$dst = $ws.Range("A626:B724")
$src.Copy($dst)
$ws.Rows("1011:1109").Delete()

# --- 2. Refresh the AutoFilter range to cover the full data extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:B1010").AutoFilter()

# --- 3. Update the hidden _FilterDatabase defined name to match.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$B`$1010"
    }
}

# --- 4. Point the PivotTable's source range at the new extent.
$ws1 = $wb.Worksheets.Item("Planilha1")
$pt = $ws1.PivotTables().Item(1)
$pt.SourceData = "Sheet1!A1:B1010"

# --- 5. Update the visible window scroll position + selection to match
#        the author's final viewport (topLeftCell A673, selection K635).
$win = $excel.ActiveWindow
$win.ScrollRow = 673
$win.ScrollColumn = 1
$ws.Range("K635").Select()
